$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Faq")

# ---- Update cell values for rows 11-13 ----
$ws.Range("B11").Value = "d"
$ws.Range("D11").Value = "tutor60@nkt.com"
$ws.Range("E11").Value = "Admin@123"
$ws.Range("F11").Value = "Pilot sess 6"

$ws.Range("D12").Value = "tutor60@nkt.com"
$ws.Range("F12").Value = "Pilot MV 6"

$ws.Range("B13").Value = "tutor"
$ws.Range("D13").Value = "tutor60@nkt.com"
$ws.Range("F13").Value = "Pilot MF 6"

# ---- Rebuild hyperlinks: the existing per-cell hyperlinks on E11/D12/E12/D13/E13
# need to be replaced by merged ranges (D12:D13 and E11:E13) plus a fresh D11 link.
# The COM bridge can't selectively delete a pre-existing hyperlink loaded from the
# file, so clear the whole collection and re-add every hyperlink (touched + untouched).
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:Admin@123") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:Admin@123") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E4"), "mailto:Admin@123") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E5"), "mailto:Admin@123") | Out-Null

$h = $ws.Hyperlinks.Add($ws.Range("D3:D6"), "mailto:tutor36@nkt.com")
$h.TextToDisplay = "tutor36@nkt.com"

$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:tutor39@nkt.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E6"), "mailto:Admin@123") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:tutor53@nkt.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D7"), "mailto:tutor53@nkt.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E7"), "mailto:Admin@123") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D8"), "mailto:tutor54@nkt.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E8"), "mailto:Admin@123") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D9"), "mailto:tutor54@nkt.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E9"), "mailto:Admin@123") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D10"), "mailto:tutor54@nkt.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E10"), "mailto:Admin@123") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D14"), "mailto:tutor54@nkt.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E14"), "mailto:Admin@123") | Out-Null

$ws.Hyperlinks.Add($ws.Range("D11"), "mailto:tutor60@nkt.com") | Out-Null

$h = $ws.Hyperlinks.Add($ws.Range("D12:D13"), "mailto:tutor60@nkt.com")
$h.TextToDisplay = "tutor60@nkt.com"

$h = $ws.Hyperlinks.Add($ws.Range("E11:E13"), "mailto:Admin@123")
$h.TextToDisplay = "Admin@123"

# ---- Restore per-cell hyperlink formatting for D/E columns that don't already
# carry it from before (D11 newly becomes a hyperlink cell like D12-D14) ----
$ws.Range("D11").Style = "Hyperlink"

# ---- Selection matches the author's last recorded cursor position ----
$ws.Range("B13").Select()
